$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 58 (hunk 0)
$ws.Range("H58").Value = 3919.2
$ws.Range("I58").Value = 111
$ws.Range("J58").Value = 4871.25
$ws.Range("K58").Value = 333
$ws.Range("L58").Value = 14613.75
$ws.Range("M58").Value = -183
$ws.Range("N58").Value = -14913.75
# Row 125 (hunk 1)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 137 (hunk 2)
$ws.Range("H137").Value = 6998.7827
$ws.Range("I137").Value = 6024.684
$ws.Range("J137").Value = 11625.75
$ws.Range("K137").Value = 18074.052
$ws.Range("L137").Value = 34877.25
$ws.Range("M137").Value = -15524.052

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 15 (hunk 3)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
# Row 16 (hunk 4)
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
# Row 36 (hunk 5)
$ws.Range("H36").Value = 7374
$ws.Range("I36").Value = 7449
$ws.Range("J36").Value = 7299
$ws.Range("K36").Value = 7449
$ws.Range("L36").Value = 7299
$ws.Range("M36").Value = -7103
$ws.Range("N36").Value = -7991
# Row 74 (hunk 6)
$ws.Range("H74").Value = 11403.2
$ws.Range("I74").Value = 5672
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 5672
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -4798
$ws.Range("N74").Value = -21748
# Row 77 (hunk 7)
$ws.Range("H77").Value = 11403.2
$ws.Range("I77").Value = 5672
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 28360
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -23992
$ws.Range("N77").Value = -108736
# Row 110 (hunk 8)
$ws.Range("H110").Value = 996.75
$ws.Range("I110").Value = 996.75
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 996.75
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1048.25
# Row 122 (hunk 9)
$ws.Range("H122").Value = 2670.3333
$ws.Range("I122").Value = 2505.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7516.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5066.5
$ws.Range("N122").Value = -13900
# Row 132 (hunk 10)
$ws.Range("H132").Value = 14104.2
$ws.Range("I132").Value = 10840.333
$ws.Range("J132").Value = 19000
$ws.Range("K132").Value = 32520.999
$ws.Range("L132").Value = 57000
$ws.Range("M132").Value = -29990.999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 12 (hunk 11)
$ws.Range("H12").Value = 434.33334
$ws.Range("I12").Value = 152.5
$ws.Range("J12").Value = 998
$ws.Range("K12").Value = 152.5
$ws.Range("L12").Value = 998
$ws.Range("M12").Value = 15.5
$ws.Range("N12").Value = -1334
# Row 33 (hunk 12)
$ws.Range("H33").Value = 20001
$ws.Range("I33").Value = 20001
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 20001
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -19665
# Row 99 (hunk 13)
$ws.Range("H99").Value = 3300
$ws.Range("I99").Value = 2350
$ws.Range("J99").Value = 4250
$ws.Range("K99").Value = 2350
$ws.Range("L99").Value = 4250
$ws.Range("M99").Value = -852
$ws.Range("N99").Value = -7246
# Row 107 (hunk 14)
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 125 (hunk 15)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 134 (hunk 16)
$ws.Range("H134").Value = 5950.4546
$ws.Range("I134").Value = 3495
$ws.Range("J134").Value = 17000
$ws.Range("K134").Value = 10485
$ws.Range("L134").Value = 51000
$ws.Range("M134").Value = -7950
$ws.Range("N134").Value = -56070

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 97 (hunk 17)
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
# Row 113 (hunk 18)
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -80
$ws.Range("N113").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 4 (hunk 19)
$ws.Range("H4").Value = 2399.8
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 2666.3333
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2666.3333
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -2890.3333
# Row 6 (hunk 20)
$ws.Range("H6").Value = 1341.6666
$ws.Range("I6").Value = 1341.6666
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1341.6666
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1228.6666
# Row 12 (hunk 21)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 13 (hunk 22)
$ws.Range("H13").Value = 223.33333
$ws.Range("I13").Value = 266.8
$ws.Range("J13").Value = 6
$ws.Range("K13").Value = 266.8
$ws.Range("L13").Value = 6
$ws.Range("M13").Value = -127.8
$ws.Range("N13").Value = -284
# Row 16 (hunk 23)
$ws.Range("H16").Value = 1341.6666
$ws.Range("I16").Value = 1341.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1341.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1091.6666
# Row 17 (hunk 24)
$ws.Range("H17").Value = 662.5
$ws.Range("I17").Value = 1100
$ws.Range("J17").Value = 225
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 225
$ws.Range("M17").Value = -932
$ws.Range("N17").Value = -561
# Row 19 (hunk 25)
$ws.Range("H19").Value = 1052
$ws.Range("I19").Value = 75
$ws.Range("J19").Value = 3006
$ws.Range("K19").Value = 75
$ws.Range("L19").Value = 3006
$ws.Range("M19").Value = 213
$ws.Range("N19").Value = -3582
# Row 20 (hunk 26)
$ws.Range("H20").Value = 25000
$ws.Range("I20").Value = 15000
$ws.Range("J20").Value = 27000
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 27000
$ws.Range("M20").Value = -14755
$ws.Range("N20").Value = -27490
# Row 23 (hunk 27)
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -77
# Row 27 (hunk 28)
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 31 (hunk 29)
$ws.Range("H31").Value = 1256.8
$ws.Range("I31").Value = 1256.8
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1256.8
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -964.8
# Row 35 (hunk 30)
$ws.Range("H35").Value = 4000
$ws.Range("I35").Value = 4000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -3702
# Row 37 (hunk 31)
$ws.Range("H37").Value = 1256.8
$ws.Range("I37").Value = 1256.8
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1256.8
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -979.8
# Row 97 (hunk 32)
$ws.Range("H97").Value = 899.5
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 999
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -503
$ws.Range("N97").Value = -1792
# Row 107 (hunk 33)
$ws.Range("H107").Value = 1433.625
$ws.Range("I107").Value = 1600
$ws.Range("J107").Value = 269
$ws.Range("K107").Value = 1600
$ws.Range("L107").Value = 269
$ws.Range("M107").Value = 320
$ws.Range("N107").Value = -4109
# Row 132 (hunk 34)
$ws.Range("H132").Value = 7166.7896
$ws.Range("I132").Value = 5582.615
$ws.Range("J132").Value = 10599.167
$ws.Range("K132").Value = 16747.845
$ws.Range("L132").Value = 31797.501
$ws.Range("M132").Value = -14217.845

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 35)
$ws.Range("H22").Value = 7624.2
$ws.Range("I22").Value = 8646
$ws.Range("J22").Value = 5726.5713
$ws.Range("K22").Value = 8646
$ws.Range("L22").Value = 5726.5713
$ws.Range("M22").Value = -8351
$ws.Range("N22").Value = -6316.5713
# Row 27 (hunk 36)
$ws.Range("H27").Value = 7624.2
$ws.Range("I27").Value = 8646
$ws.Range("J27").Value = 5726.5713
$ws.Range("K27").Value = 8646
$ws.Range("L27").Value = 5726.5713
$ws.Range("M27").Value = -8539
$ws.Range("N27").Value = -5940.5713
# Row 40 (hunk 37)
$ws.Range("H40").Value = 10829.333
$ws.Range("I40").Value = 10829.333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 10829.333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -10693.333
# Row 68 (hunk 38)
$ws.Range("H68").Value = 5466
$ws.Range("I68").Value = 5466
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5466
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4717
# Row 71 (hunk 39)
$ws.Range("H71").Value = 5466
$ws.Range("I71").Value = 5466
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 27330
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -23586

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 15 (hunk 40)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 122 (hunk 41)
$ws.Range("H122").Value = 12004
$ws.Range("I122").Value = 12004
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 36012
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -33562
$ws.Range("N122").ClearContents()
